# Insert a new data row at row 45 (weekly update: a new price observation
# was recorded, pushing the existing rows 45-77 down to 46-78).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 45..77 down to 46..78 by inserting a blank row at 45.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new observation.
$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 45271
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 300000000
$ws.Range("G45").Value = "Espárragos"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 900
$ws.Range("K45").Value = 1300
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = 1456
$ws.Range("N45").Value = "$/kilo"
$ws.Range("O45").Value = "Región de Ñuble"
$ws.Range("P45").Value = 1456
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"
